# Daten aktualisiert am 2024-01-31
# Append three new ticker rows to the end of the existing list on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$newTickers = @("IMX-USD", "TAO-USD", "MNT-USD")

for ($i = 0; $i -lt $newTickers.Length; $i++) {
    $row = $lastRow + $i + 1
    $ws.Cells.Item($row, 1).Value = $newTickers[$i]
}
